$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.393.79"
$ws.Range("E2").Value = "  -6.27%  "

$ws.Range("D3").Value = "2.890.08"
$ws.Range("E3").Value = "  -4.14%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.70"
$ws.Range("E5").Value = "  -3.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "122.81"
$ws.Range("E6").Value = "  -4.71%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").Value = "2.884.56"
$ws.Range("E8").Value = "  -4.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  -0.60%  "

$ws.Range("E10").Value = "  -7.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.75"
$ws.Range("E11").Value = "  -8.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.438"
$ws.Range("E12").Value = "  +1.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000212"
$ws.Range("E13").Value = "  -4.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.28"
$ws.Range("E14").Value = "  -1.71%  "

$ws.Range("E15").Value = "  +1.26%  "

$ws.Range("D16").Value = "3.384.60"
$ws.Range("E16").Value = "  -3.53%  "

$ws.Range("D17").Value = "2.902.14"
$ws.Range("E17").Value = "  -3.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.58"
$ws.Range("E18").Value = "  +5.75%  "

$ws.Range("D19").Value = "57.437.70"
$ws.Range("E19").Value = "  -6.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "408.69"
$ws.Range("E20").Value = "  -6.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.90"
$ws.Range("E21").Value = "  -2.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.670"
$ws.Range("E22").Value = "  +1.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.84"
$ws.Range("E23").Value = "  -4.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.87"
$ws.Range("E24").Value = "  +2.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "77.09"
$ws.Range("E25").Value = "  -2.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.45"
$ws.Range("E28").Value = "  -2.08%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.93"
$ws.Range("E29").Value = "  +3.00%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.22"
$ws.Range("E30").Value = "  +0.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.03"
$ws.Range("E31").Value = "  -3.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.68"
$ws.Range("E32").Value = "  -3.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0978"
$ws.Range("E33").Value = "  +3.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.910"
$ws.Range("E34").Value = "  -4.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.37"
$ws.Range("E35").Value = "  -3.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.00"
$ws.Range("E36").Value = "  -11.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.19"
$ws.Range("E37").Value = "  -3.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.45"
$ws.Range("E38").Value = "  +9.30%  "

$ws.Range("D39").Value = "0.0₃0620"
$ws.Range("E39").Value = "  -8.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0343"
$ws.Range("E40").Value = "  -5.65%  "

$ws.Range("E41").Value = "  -2.01%  "

$ws.Range("D42").Value = "2.622.79"
$ws.Range("E42").Value = "  -0.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "360.19"
$ws.Range("E43").Value = "  -3.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.39"
$ws.Range("E44").Value = "  -1.57%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "119.93"
$ws.Range("E46").Value = "  +0.37%  "

$ws.Range("E47").Value = "  -3.15%  "

$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.93"
$ws.Range("E49").Value = "  -1.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.79"
$ws.Range("E50").Value = "  -3.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.95"
$ws.Range("E51").Value = "  -3.83%  "
